$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B264").Value = 48719
$ws.Range("E264").Value = 353.35
$ws.Range("F264").Value = -81
$ws.Range("G264").Value = -23955.75
$ws.Range("B265").Value = 64979
$ws.Range("E265").Value = 314.41
$ws.Range("F265").Value = 82
$ws.Range("G265").Value = 24251.5
$ws.Range("B313").Value = 62997
$ws.Range("F313").Value = 72
$ws.Range("G313").Value = 22020.48
$ws.Range("B314").Value = 57854
$ws.Range("F314").Value = 2
$ws.Range("G314").Value = 611.6799999999999
$ws.Range("B316").Value = 57077
$ws.Range("D316").Value = 93.08
$ws.Range("E316").Value = 111.2
$ws.Range("F316").Value = 1
$ws.Range("G316").Value = 93.08
$ws.Range("B317").Value = 61610
$ws.Range("D317").Value = 102.71
$ws.Range("E317").Value = 122.71
$ws.Range("F317").Value = -58
$ws.Range("G317").Value = -5957.18
$ws.Range("B318").Value = 63565
$ws.Range("E318").Value = 109.19
$ws.Range("F318").Value = 60
$ws.Range("G318").Value = 6162.6
$ws.Range("B346").Value = 63520
$ws.Range("E346").Value = 153.4
$ws.Range("F346").Value = 97
$ws.Range("G346").Value = 13995.16
$ws.Range("B347").Value = 55373
$ws.Range("E347").Value = 163.62
$ws.Range("F347").Value = -94
$ws.Range("G347").Value = -13562.32
$ws.Range("B350").Value = 63531
$ws.Range("F350").Value = 80
$ws.Range("G350").Value = 11478.4
$ws.Range("B352").Value = 63571
$ws.Range("F352").Value = 29
$ws.Range("G352").Value = 4160.92
$ws.Range("B355").Value = 55356
$ws.Range("E355").Value = 54.04
$ws.Range("F355").Value = -158
$ws.Range("G355").Value = -7527.12
$ws.Range("B356").Value = 63510
$ws.Range("E356").Value = 50.66
$ws.Range("F356").Value = 167
$ws.Range("G356").Value = 7955.88
$ws.Range("B372").Value = 63652
$ws.Range("E372").Value = 55.42
$ws.Range("F372").Value = 250
$ws.Range("G372").Value = 13032.5
$ws.Range("B373").Value = 57885
$ws.Range("E373").Value = 62.28
$ws.Range("F373").Value = 4
$ws.Range("G373").Value = 208.52
$ws.Range("B379").Value = 61608
$ws.Range("E379").Value = 154.12
$ws.Range("F379").Value = -56
$ws.Range("G379").Value = -7224.56
$ws.Range("B380").Value = 63564
$ws.Range("E380").Value = 137.16
$ws.Range("F380").Value = 57
$ws.Range("G380").Value = 7353.57
$ws.Range("B389").Value = 62865
$ws.Range("F389").Value = 151
$ws.Range("G389").Value = 12051.31
$ws.Range("B390").Value = 57817
$ws.Range("F390").Value = 3
$ws.Range("G390").Value = 239.43
$ws.Range("B419").Value = 63007
$ws.Range("F419").Value = 984
$ws.Range("G419").Value = 168588.72
$ws.Range("B420").Value = 57856
$ws.Range("F420").Value = 2
$ws.Range("G420").Value = 342.66
$ws.Range("B421").Value = 63008
$ws.Range("F421").Value = 504
$ws.Range("G421").Value = 76189.67999999999
$ws.Range("B422").Value = 57857
$ws.Range("F422").Value = 3
$ws.Range("G422").Value = 453.51
$ws.Range("B431").Value = 53082
$ws.Range("C431").Value = "HUL-VIM BAR MULTIPACK FW 4X200G"
$ws.Range("F431").Value = 1
$ws.Range("G431").Value = 59.47
$ws.Range("B432").Value = 63102
$ws.Range("C432").Value = "HUL-Vim Bar Multipack Fw 4X200G"
$ws.Range("F432").Value = 36
$ws.Range("G432").Value = 2140.92
$ws.Range("B583").Value = 53263
$ws.Range("E583").Value = 15.29
$ws.Range("F583").Value = -309
$ws.Range("G583").Value = -3958.29
$ws.Range("B584").Value = 65066
$ws.Range("E584").Value = 13.61
$ws.Range("F584").Value = 313
$ws.Range("G584").Value = 4009.53
$ws.Range("B590").Value = 64922
$ws.Range("E590").Value = 20.98
$ws.Range("F590").Value = 207
$ws.Range("G590").Value = 4084.11
$ws.Range("B591").Value = 45706
$ws.Range("E591").Value = 23.58
$ws.Range("F591").Value = -202
$ws.Range("G591").Value = -3985.46
$ws.Range("B599").Value = 45709
$ws.Range("E599").Value = 15.69
$ws.Range("F599").Value = -300
$ws.Range("G599").Value = -3945
$ws.Range("B600").Value = 64925
$ws.Range("E600").Value = 13.97
$ws.Range("F600").Value = 302
$ws.Range("G600").Value = 3971.3
$ws.Range("B601").Value = 64919
$ws.Range("E601").Value = 27.97
$ws.Range("F601").Value = 224
$ws.Range("G601").Value = 5891.2
$ws.Range("B602").Value = 45702
$ws.Range("E602").Value = 31.43
$ws.Range("F602").Value = -215
$ws.Range("G602").Value = -5654.5
$ws.Range("B709").Value = 64833
$ws.Range("E709").Value = 34.9
$ws.Range("F709").Value = 99
$ws.Range("G709").Value = 3250.17
$ws.Range("B710").Value = 60025
$ws.Range("E710").Value = 37.22
$ws.Range("F710").Value = -98
$ws.Range("G710").Value = -3217.34
$ws.Range("B715").Value = 60031
$ws.Range("E715").Value = 111.69
$ws.Range("F715").Value = -5
$ws.Range("G715").Value = -492.5
$ws.Range("B716").Value = 64836
$ws.Range("E716").Value = 104.71
$ws.Range("F716").Value = 7
$ws.Range("G716").Value = 689.5
$ws.Range("B872").Value = 65362
$ws.Range("F872").Value = 2
$ws.Range("G872").Value = 81.73999999999999
$ws.Range("B873").Value = 65079
$ws.Range("F873").Value = 21
$ws.Range("G873").Value = 858.27
